# "This is my Ninth day first commit"
# Append the day-9 rows (85-95) to the TaskList sheet, continuing the
# existing task log table (S.No, Date, Task, Link to video, Link to
# material, Time Taken, Errors Y/N, Error report link).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TaskList")

# Data for the eleven new rows (84 - 94 in S.No, sheet rows 85 - 95),
# listed in the order the author actually typed the "Task" cells (the
# shared-string table records new strings in first-use order, and the
# author filled row 90's task in before going back and filling row 89's
# - e.g. typed the two Register/FrontEnd tasks out of row order).
# Columns: SNo, Date, Task, Material link, Time taken
$rows = @(
  @{ R=85; SNo=84; Task="Modification of User Entity Class"; Link=""; Time="10Minutes" },
  @{ R=86; SNo=85; Task="Modification of UserDAOImpl class"; Link=""; Time="15Minutes" },
  @{ R=87; SNo=86; Task="Create UserTestCase for inset, update, delete, retrive"; Link=""; Time="40Minutes" },
  @{ R=88; SNo=87; Task="Product Testcase for update and delete"; Link=""; Time="15Minutes" },
  @{ R=90; SNo=89; Task="Run the Register page to register as supplier and customer"; Link=""; Time="15Minutes" },
  @{ R=89; SNo=88; Task="Create FronEndUserController class"; Link=""; Time="90Minutes" },
  @{ R=91; SNo=90; Task="Configuration of Webflow"; Link="http://www.jcombat.com/spring/create-your-first-spring-web-flow-based-web-application"; Time="30Minutes" },
  @{ R=92; SNo=91; Task="Designed billing, preRegister, navbar, welcome pages"; Link=""; Time="60Minutes" },
  @{ R=93; SNo=92; Task="Created membershipflow.xml file"; Link="http://www.jcombat.com/spring/create-your-first-spring-web-flow-based-web-application"; Time="30Minutes" },
  @{ R=94; SNo=93; Task="Created Address model"; Link=""; Time="15Minutes" },
  @{ R=95; SNo=94; Task="Created RegisterModel class"; Link=""; Time="10Minutes" }
)

# Rows whose wrapped Task text spans two lines in the real workbook
# (taller row height, same as similarly long tasks elsewhere in the sheet).
$tallRows = @(87, 90, 91, 92, 93)

$entryDate = 42801

foreach ($row in $rows) {
  $r = $row.R

  # Copy the formatting of an existing data row so borders / fonts /
  # number formats / wrap settings match the rest of the table. The
  # "S.No" column alternates between two border styles on every other
  # row (matches rows 83/84 and continues down), so pick whichever
  # template row lines up with the new row's position in that cadence.
  if (($r % 2) -eq 1) {
    $ws.Range("A83:H83").Copy()
  } else {
    $ws.Range("A84:H84").Copy()
  }
  $ws.Range("A" + $r + ":H" + $r).PasteSpecial(-4122)  # xlPasteFormats

  $ws.Cells.Item($r, 1).Value = $row.SNo
  $ws.Cells.Item($r, 2).Value = $entryDate
  $ws.Cells.Item($r, 3).Value = $row.Task
  if ($row.Link -ne "") {
    $ws.Cells.Item($r, 5).Value = $row.Link
  }
  $ws.Cells.Item($r, 6).Value = $row.Time

  if ($tallRows -contains $r) {
    $ws.Rows.Item($r).RowHeight = 28.8
  }
}

$excel.CutCopyMode = $false

# Restore the on-screen selection to where the author ended up after
# typing the new rows.
$ws.Activate() | Out-Null
$ws.Range("D91").Select() | Out-Null
